$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 194983.44
$ws.Range("J17").Value = 200861.73
$ws.Range("L17").Value = 602585.1900000001
$ws.Range("N17").Value = -602921.1900000001
$ws.Range("H19").Value = 1500
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 2000
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = -825
$ws.Range("N19").Value = -2350
$ws.Range("H33").Value = 3490012.5
$ws.Range("I33").Value = 4590874.5
$ws.Range("K33").Value = 4590874.5
$ws.Range("M33").Value = -4590645.5
$ws.Range("H92").Value = 1953587.4
$ws.Range("I92").Value = 744485.9
$ws.Range("K92").Value = 744485.9
$ws.Range("M92").Value = -743237.9
$ws.Range("H132").Value = 1605
$ws.Range("I132").Value = 1440.9412
$ws.Range("K132").Value = 4322.8236
$ws.Range("M132").Value = -1792.8236
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 79179.30499999999
$ws.Range("I45").Value = 93030.09
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 93030.09
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -92653.09
$ws.Range("N45").Value = -3754
$ws.Range("H61").Value = 899950
$ws.Range("I61").Value = 24020.283
$ws.Range("K61").Value = 24020.283
$ws.Range("M61").Value = -23808.283
$ws.Range("H88").Value = 2924.25
$ws.Range("J88").Value = 2924.25
$ws.Range("L88").Value = 2924.25
$ws.Range("N88").Value = -3736.25
$ws.Range("H91").Value = 2924.25
$ws.Range("J91").Value = 2924.25
$ws.Range("L91").Value = 2924.25
$ws.Range("N91").Value = -5732.25
$ws.Range("H122").Value = 2235.3076
$ws.Range("I122").Value = 2244.76
$ws.Range("K122").Value = 6734.280000000001
$ws.Range("M122").Value = -4284.280000000001
$ws.Range("H136").Value = 899950
$ws.Range("I136").Value = 24020.283
$ws.Range("K136").Value = 72060.849
$ws.Range("M136").Value = -69510.849
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 3250
$ws.Range("J39").Value = 6000
$ws.Range("L39").Value = 6000
$ws.Range("N39").Value = -6778
$ws.Range("H94").Value = 2060.4375
$ws.Range("I94").Value = 1113.8334
$ws.Range("K94").Value = 1113.8334
$ws.Range("M94").Value = -662.8334
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 20000000
$ws.Range("J4").Value = 20000000
$ws.Range("L4").Value = 20000000
$ws.Range("N4").Value = -20000224
$ws.Range("H16").Value = 15876072
$ws.Range("I16").Value = 35715784
$ws.Range("J16").Value = 4302
$ws.Range("K16").Value = 35715784
$ws.Range("L16").Value = 4302
$ws.Range("M16").Value = -35715497
$ws.Range("N16").Value = -4876
$ws.Range("H22").Value = 573.5484
$ws.Range("J22").Value = 650.8
$ws.Range("L22").Value = 650.8
$ws.Range("N22").Value = -1350.8
$ws.Range("H94").Value = 657.1
$ws.Range("I94").Value = 879.3333
$ws.Range("K94").Value = 879.3333
$ws.Range("M94").Value = -428.3333
$ws.Range("H99").Value = 57779268
$ws.Range("I99").Value = 3335068.8
$ws.Range("K99").Value = 3335068.8
$ws.Range("M99").Value = -3333570.8
$ws.Range("H107").Value = 1559.5555
$ws.Range("I107").Value = 1559.5555
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1559.5555
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 360.4445000000001
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 15876072
$ws.Range("I113").Value = 35715784
$ws.Range("J113").Value = 4302
$ws.Range("K113").Value = 35715784
$ws.Range("L113").Value = 4302
$ws.Range("M113").Value = -35713614
$ws.Range("N113").Value = -8642
$ws.Range("H122").Value = 1872.6
$ws.Range("I122").Value = 1872.6
$ws.Range("K122").Value = 5617.799999999999
$ws.Range("M122").Value = -3167.799999999999
$ws.Range("H126").Value = 57779268
$ws.Range("I126").Value = 3335068.8
$ws.Range("K126").Value = 10005206.4
$ws.Range("M126").Value = -10002736.4
$ws.Range("H132").Value = 3110.0557
$ws.Range("I132").Value = 2934.3333
$ws.Range("J132").Value = 3285.7778
$ws.Range("K132").Value = 8802.999899999999
$ws.Range("L132").Value = 9857.3334
$ws.Range("M132").Value = -6272.999899999999
$ws.Range("N132").Value = -14917.3334
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 503.4
$ws.Range("I14").Value = 503.4
$ws.Range("K14").Value = 1510.2
$ws.Range("M14").Value = -1337.2
$ws.Range("H24").Value = 2055.15
$ws.Range("J24").Value = 2097.5789
$ws.Range("L24").Value = 6292.736699999999
$ws.Range("N24").Value = -6752.736699999999
$ws.Range("H112").Value = 5669.778
$ws.Range("I112").Value = 8766.666999999999
$ws.Range("K112").Value = 26300.001
$ws.Range("M112").Value = -25192.001
$ws.Range("H134").Value = 2256.55
$ws.Range("I134").Value = 1743.7368
$ws.Range("K134").Value = 5231.2104
$ws.Range("M134").Value = -161.2103999999999
$ws.Range("H136").Value = 9120.321
$ws.Range("I136").Value = 5210.875
$ws.Range("J136").Value = 14332.917
$ws.Range("K136").Value = 15632.625
$ws.Range("L136").Value = 42998.751
$ws.Range("M136").Value = -10532.625
$ws.Range("N136").Value = -53198.751
$ws.Range("H140").Value = 1979.6666
$ws.Range("J140").Value = 3124.75
$ws.Range("L140").Value = 9374.25
$ws.Range("N140").Value = -19734.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 3852.9
$ws.Range("I41").Value = 2614.3333
$ws.Range("K41").Value = 2614.3333
$ws.Range("M41").Value = -2259.3333
$ws.Range("H57").Value = 76666.336
$ws.Range("J57").Value = 99999
$ws.Range("L57").Value = 99999
$ws.Range("N57").Value = -101639
$ws.Range("H97").Value = 84025.96000000001
$ws.Range("I97").Value = 63027.688
$ws.Range("K97").Value = 63027.688
$ws.Range("M97").Value = -62531.688
$ws.Range("H99").Value = 9369.538
$ws.Range("I99").Value = 7233.75
$ws.Range("K99").Value = 7233.75
$ws.Range("M99").Value = -4987.75
$ws.Range("H113").Value = 1255
$ws.Range("I113").Value = 891.6667
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 891.6667
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 1278.3333
$ws.Range("N113").Value = -6140
$ws.Range("H122").Value = 1183.4706
$ws.Range("I122").Value = 940.73334
$ws.Range("K122").Value = 2822.20002
$ws.Range("M122").Value = -372.2000200000002
$ws.Range("H136").Value = 38997.668
$ws.Range("J136").Value = 38997.668
$ws.Range("L136").Value = 116993.004
$ws.Range("N136").Value = -122093.004
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4283.6787
$ws.Range("I22").Value = 1265.6
$ws.Range("K22").Value = 1265.6
$ws.Range("M22").Value = -970.5999999999999
$ws.Range("H27").Value = 4283.6787
$ws.Range("I27").Value = 1265.6
$ws.Range("K27").Value = 1265.6
$ws.Range("M27").Value = -1158.6
$ws.Range("H43").Value = 1183749.5
$ws.Range("I43").Value = 200000
$ws.Range("J43").Value = 1442631
$ws.Range("K43").Value = 200000
$ws.Range("L43").Value = 1442631
$ws.Range("M43").Value = -199807
$ws.Range("N43").Value = -1443017
$ws.Range("H61").Value = 1410.9546
$ws.Range("I61").Value = 1397
$ws.Range("K61").Value = 1397
$ws.Range("M61").Value = -1195
$ws.Range("H113").Value = 1410.9546
$ws.Range("I113").Value = 1397
$ws.Range("K113").Value = 1397
$ws.Range("M113").Value = 773
$ws.Range("H136").Value = 58803
$ws.Range("I136").Value = 73961.71000000001
$ws.Range("J136").Value = 5747.5
$ws.Range("K136").Value = 221885.13
$ws.Range("L136").Value = 17242.5
$ws.Range("M136").Value = -219335.13
$ws.Range("N136").Value = -22342.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2431.1738
$ws.Range("I122").Value = 2460.85
$ws.Range("K122").Value = 7382.549999999999
$ws.Range("M122").Value = -4932.549999999999
$ws.Range("H125").Value = 89139
$ws.Range("J125").Value = 89139
$ws.Range("L125").Value = 89139
$ws.Range("N125").Value = -98979
$ws.Range("H132").Value = 2448.04
$ws.Range("I132").Value = 2080.2222
$ws.Range("J132").Value = 3393.8572
$ws.Range("K132").Value = 6240.6666
$ws.Range("L132").Value = 10181.5716
$ws.Range("M132").Value = -3710.6666
$ws.Range("N132").Value = -15241.5716
